# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect that the
# de-de handback has completed ("Generate Report for Handback"):
#   - Overview!Status (zh-cn / de-de columns) now reads
#     "Handed back: in sync with en-US" instead of "Ready for handoff"
#   - Both the zh-cn and de-de detail sheets get their "Latest Target File"
#     and "Latest Handback File" columns filled in (with a hyperlink on the
#     target-file cell), and their "Latest Handback DateTime" updated.
#   - A few columns are widened so the new values are easier to read.

$wb = $excel.ActiveWorkbook

# Helper: Excel quantizes ColumnWidth (character units) to whole pixels
# before storing it back out as the OOXML <col width> value (pixels = the
# rounded character width * max-digit-width, then width = (pixels+5)/mdw).
# Given a *desired* stored width, compute the ColumnWidth to assign that
# lands on the closest achievable pixel bucket.
function Set-ColWidth {
    param($ws, $colIndex, $desiredWidth)

    $mdw = 6
    $px = [Math]::Round(($desiredWidth * $mdw) - 5)
    if ($px -lt 0) { $px = 0 }
    $lo = ($px - 0.5) / $mdw
    $hi = ($px + 0.5) / $mdw
    $chars = ($lo + $hi) / 2
    $ws.Columns.Item($colIndex).ColumnWidth = $chars
}

$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/263f9135c1094ede58d26b82153c343b78096741/e2e/acc6e68f-631e-43f6-a6c7-ea97c063043d.md"
$mdDisplay = "acc6e68f-631e-43f6-a6c7-ea97c063043d.md"

# ---------------------------------------------------------------------
# Overview sheet: Status column for zh-cn (E) and de-de (F) now reports
# that the handback is in sync with en-US.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

Set-ColWidth $wsOverview 5 29.9777047293527
Set-ColWidth $wsOverview 6 29.9777047293527

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Latest Target File (I2): link back to the source markdown file.
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl, "", "", $mdDisplay)

# Latest Handback File (J2): the generated zh-cn xlf that was handed back.
$wsZhCn.Range("J2").Value = "acc6e68f-631e-43f6-a6c7-ea97c063043d.2f8f1c07ecb678d0e488e04d22bd1c1d9a7fe724.zh-cn.xlf"

# Latest Handback DateTime (K2)
$wsZhCn.Range("K2").Value = "2016-08-20 15:01:02"

Set-ColWidth $wsZhCn 3 29.9777047293527
Set-ColWidth $wsZhCn 9 40
Set-ColWidth $wsZhCn 10 40

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Latest Target File (I2): link back to the source markdown file.
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl, "", "", $mdDisplay)

# Latest Handback File (J2): the generated de-de xlf that was handed back.
$wsDeDe.Range("J2").Value = "acc6e68f-631e-43f6-a6c7-ea97c063043d.2f8f1c07ecb678d0e488e04d22bd1c1d9a7fe724.de-de.xlf"

# Latest Handback DateTime (K2)
$wsDeDe.Range("K2").Value = "2016-08-20 15:01:11"

Set-ColWidth $wsDeDe 3 29.9777047293527
Set-ColWidth $wsDeDe 9 40
Set-ColWidth $wsDeDe 10 40

$wb.Save()
